$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.883.64"
$ws.Range("E2").Value = "  -3.09%  "
$ws.Range("D3").Value = "3.324.63"
$ws.Range("E3").Value = "  -5.30%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "551.97"
$ws.Range("E5").Value = "  -4.56%  "
$ws.Range("D6").Value = "172.10"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("E7").Value = "  -3.75%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "3.317.23"
$ws.Range("E9").Value = "  -5.35%  "
$ws.Range("D10").Value = "0.621"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").Value = "52.95"
$ws.Range("E12").Value = "  -4.51%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "9.01"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "3.853.59"
$ws.Range("E15").Value = "  -5.47%  "
$ws.Range("D16").Value = "18.19"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("E17").Value = "  -3.56%  "
$ws.Range("D18").Value = "3.329.94"
$ws.Range("E18").Value = "  -5.01%  "
$ws.Range("D19").Value = "11.80"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "63.814.08"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "0.969"
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("D22").Value = "428.43"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("D23").Value = "4.58"
$ws.Range("E23").Value = "  +6.45%  "
$ws.Range("D24").Value = "4.08"
$ws.Range("E24").Value = "  -4.19%  "
$ws.Range("D25").Value = "84.15"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "13.41"
$ws.Range("E26").Value = "  +2.81%  "
$ws.Range("D27").Value = "10.60"
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("D28").Value = "2.81"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "8.57"
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("D30").Value = "29.60"
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("D31").Value = "6.60"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("D32").Value = "592.10"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("D33").Value = "11.40"
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").Value = "0.106"
$ws.Range("D35").Value = "58.14"
$ws.Range("E35").Value = "  -2.52%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "0.142"
$ws.Range("E37").Value = "  -9.29%  "
$ws.Range("D38").Value = "35.23"
$ws.Range("E38").Value = "  -5.55%  "
$ws.Range("D39").Value = "0.0₃0747"
$ws.Range("E39").Value = "  -6.42%  "
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("D41").Value = "0.363"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("D42").Value = "3.094.09"
$ws.Range("E42").Value = "  -6.35%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "2.77"
$ws.Range("E44").Value = "  -5.09%  "
$ws.Range("D45").Value = "0.0405"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("D46").Value = "3.16"
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").Value = "2.57"
$ws.Range("E49").Value = "  -5.49%  "
$ws.Range("D50").Value = "8.15"
$ws.Range("E50").Value = "  -4.91%  "
$ws.Range("D51").Value = "132.41"
$ws.Range("E51").Value = "  -5.57%  "
